$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 8 and row 9: A, M, Q, R, Z, AB, AC
$cols = @("A", "M", "Q", "R", "Z", "AB", "AC")

foreach ($col in $cols) {
    $cell8 = $ws.Range($col + "8")
    $cell9 = $ws.Range($col + "9")

    $v8 = $cell8.Value2
    $v9 = $cell9.Value2

    $cell8.Value2 = $v9
    $cell9.Value2 = $v8
}
